$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refactor: "QuestionModel" class box becomes "AnswerModel" ---
$ws.Range("B2").Value = "AnswerModel"
$ws.Range("J3").Value = "Answer Model"

# --- New attribute on the CharacterModel / row 6 of the diagram ---
$ws.Range("B6").Value = "userAnswer"

# --- Apply the "Good" (green) cell style to the task cells that now
#     belong to the shown/active task list (J5, J6, J9) ---
$ws.Range("J5").Style = "Good"
$ws.Range("J6").Style = "Good"
$ws.Range("J9").Style = "Good"

# --- New list of items describing the "answer options as buttons" work ---
$ws.Range("A14").Value = "getCharacters"
$ws.Range("A15").Value = "generateQuestions"
$ws.Range("A16").Value = "quizController"
$ws.Range("A17").Value = "displayAnswer"

# --- Widen column J to fit the longer task descriptions ---
$ws.Columns.Item(10).ColumnWidth = 24.41666666666667

# --- Move the active selection down below the new content ---
$ws.Range("A20").Select()
